# Edit script for DCI_Connectivity_Workflow_TUBRLC.docx
$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: merge the three runs "____ 1. " + "Select" + " watercourse"
# into a single run "____ 1. Select watercourse" (bold formatting kept).
# ---------------------------------------------------------------
$find = $d.Content
$find.Find.ClearFormatting()
$find.Find.Execute("____ 1. Select watercourse", $false, $false, $false, $false, $false, $true, 1, $false, "____ 1. Select watercourse", 2) | Out-Null

# A paragraph further up the document that already carries the
# "NOTE: ..." bullet-list formatting we want to reuse (pStyle
# "List Paragraph" + numId 1) so new NOTE paragraphs join the same list.
$noteTemplate = $d.Paragraphs(13)

# ---------------------------------------------------------------
# Change 2: add a new NOTE bullet after step 6 (Calculate stream length)
# ---------------------------------------------------------------
$step6 = $d.Content
$step6.Find.ClearFormatting()
$step6.Find.Execute("____ 6. Calculate stream length", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$step6Para = $step6.Paragraphs(1)

$step6Para.Range.InsertParagraphAfter() | Out-Null
$newNote1 = $step6Para.Next()
$newNote1.Style = "List Paragraph"
$newNote1.Range.ListFormat.ApplyListTemplateWithLevel($noteTemplate.Range.ListFormat.ListTemplate, $true, 1, $false, $false) | Out-Null
$newNote1.Range.Text = "NOTE: If network contains lakes/reservoir flowlines, set their stream length to 0 km."
$newNote1.Range.Font.Size = 10
$newNote1.Range.Font.SizeBi = 10

# ---------------------------------------------------------------
# Change 3: trim the trailing sentence from step 7 (Add To/From Node
# Fields) and add a new NOTE bullet with the replacement guidance.
# ---------------------------------------------------------------
$trim = $d.Content
$trim.Find.ClearFormatting()
$oldSentence = ". Note that the terminus should have a ‘To_Node’ that is one greater than the number of links in the network."
$trim.Find.Execute($oldSentence, $false, $false, $false, $false, $false, $true, 1, $false, ". ", 2) | Out-Null

$step7 = $d.Content
$step7.Find.ClearFormatting()
$step7.Find.Execute("____ 7. Add To/From Node Fields", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$step7Para = $step7.Paragraphs(1)

$step7Para.Range.InsertParagraphAfter() | Out-Null
$newNote2 = $step7Para.Next()
$newNote2.Style = "List Paragraph"
$newNote2.Range.ListFormat.ApplyListTemplateWithLevel($noteTemplate.Range.ListFormat.ListTemplate, $true, 1, $false, $false) | Out-Null
$newNote2.Range.Text = "NOTE: The ‘To_Node’ at the terminus is the only node that shouldn’t appear in the ‘From_Node’ list."
$newNote2.Range.Font.Size = 10
$newNote2.Range.Font.SizeBi = 10

# ---------------------------------------------------------------
# Change 4: drop the stale lastRenderedPageBreak cached before step 8
# (re-writing the run text regenerates it without the cached break).
# ---------------------------------------------------------------
$step8 = $d.Content
$step8.Find.ClearFormatting()
$step8.Find.Execute("____ 8. Convert watercourse to point layer", $false, $false, $false, $false, $false, $true, 1, $false, "____ 8. Convert watercourse to point layer", 2) | Out-Null

Write-Host "Edits applied"
